$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2.437664035705469,1.700078601593715,1.456920224378268,1.346585525682575),
    @(2.444871584426221,1.703624045046675,1.458610037681003,1.347817107110185),
    @(2.41651902216063,1.689118343792704,1.451047547992394,1.342443321396813),
    @(2.446396964240126,1.705623931229836,1.446216490933229,1.350032931609149),
    @(2.453820867654344,1.709811521561945,1.448875947771873,1.352104465373867),
    @(2.445607710036875,1.705657750267797,1.461064592524589,1.3502050541849),
    @(2.428272206215448,1.691569675101092,1.463279494679763,1.34067164994444),
    @(2.438785451754285,1.699117520271848,1.455019439359377,1.344961376309288),
    @(2.350206967072472,1.664463245616165,1.44357601379494,1.338037341956309),
    @(2.355463786409341,1.670736457539431,1.446995113239449,1.339872308311112),
    @(2.232917583795264,1.617529040115271,1.42465114317971,1.324747553747766),
    @(2.338634367234772,1.658285052561531,1.440235713151947,1.335620606451476)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("B$row").Value = $data[$i][0]
    $ws.Range("C$row").Value = $data[$i][1]
    $ws.Range("D$row").Value = $data[$i][2]
    $ws.Range("E$row").Value = $data[$i][3]
}
